$d = $word.ActiveDocument

$pairs = @(
    @("2024-10-14 Monday", "2024-10-15 Tuesday"),
    @("277×7=1939", "485×5=2425"),
    @("855×8=6840", "315×4=1260"),
    @("812×4=3248", "383×3=1149"),
    @("347×7=2429", "255×2=510"),
    @("964×5=4820", "296×5=1480"),
    @("932×8=7456", "578×9=5202"),
    @("144×8=1152", "981×8=7848"),
    @("716×5=3580", "999×4=3996"),
    @("235×8=1880", "628×7=4396"),
    @("564×7=3948", "928×3=2784"),
    @("355×8=2840", "972×2=1944"),
    @("449×9=4041", "657×6=3942"),
    @("462×9=4158", "267×5=1335"),
    @("457×3=1371", "864×8=6912"),
    @("256×9=2304", "444×3=1332"),
    @("342×8=2736", "209×4=836"),
    @("963×5=4815", "883×7=6181"),
    @("389×9=3501", "694×6=4164"),
    @("691×3=2073", "655×2=1310"),
    @("609×8=4872", "695×5=3475"),
    @("334×8=2672", "259×7=1813"),
    @("702×3=2106", "535×2=1070"),
    @("139×6=834", "712×4=2848"),
    @("445×7=3115", "381×6=2286"),
    @("958×3=2874", "381×6=2286")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
